$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert the new "apoio_*" columns (std/min/max) right after the
#    existing "apoio_medio" column (L). This shifts the old M..P
#    (contribuicoes, media_contribuicoes, menor_ano, maior_ano)
#    to P..S, and - because Excel copies the format of the column
#    being displaced - the 3 new cells inherit style "3" (the same
#    R$ #,##0.00 style used by column L) automatically.
# ------------------------------------------------------------------
$ws.Range("M1:O1").EntireColumn.Insert()

# ------------------------------------------------------------------
# 2. Insert the new "contribuicoes_*" columns (std/min/max) right
#    after "contribuicoes_med" (old "media_contribuicoes", now Q).
#    This shifts old "menor_ano"/"maior_ano" (now R,S) further right
#    to U,V, and the 3 new cells inherit style "1" (#,##0) from Q.
# ------------------------------------------------------------------
$ws.Range("R1:T1").EntireColumn.Insert()

# ------------------------------------------------------------------
# 3. Header row text (renames + new headers)
# ------------------------------------------------------------------
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"
# L1 "apoio_medio" unchanged
$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"
# P1 "contribuicoes" unchanged
$ws.Range("Q1").Value = "contribuicoes_med"
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"
# U1 "menor_ano" / V1 "maior_ano" unchanged

# ------------------------------------------------------------------
# 4. Data rows: B column author-classification label unchanged text
#    (values already correct after the shift - only new columns &
#    the recomputed "apoio_medio" column need new values).
# ------------------------------------------------------------------

# Row 2 (coletivo)
$ws.Range("L2").Value = 88.06007413874362
$ws.Range("M2").Value = 47.96955813517068
$ws.Range("N2").Value = 14.90596347946683
$ws.Range("O2").Value = 254.2443749773306
$ws.Range("R2").Value = 325.0284071787353
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 2015

# Row 3 (empresa)
$ws.Range("L3").Value = 89.8194628153171
$ws.Range("M3").Value = 37.89772641853159
$ws.Range("N3").Value = 16.18065842403185
$ws.Range("O3").Value = 233.3973531230909
$ws.Range("R3").Value = 486.1905468340719
$ws.Range("S3").Value = 1
$ws.Range("T3").Value = 7954

# Row 4 (feminino)
$ws.Range("L4").Value = 67.58055662882595
$ws.Range("M4").Value = 23.28873679351738
$ws.Range("N4").Value = 18.47818326605706
$ws.Range("O4").Value = 154.8484188303038
$ws.Range("R4").Value = 95.76839418448722
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 453

# Row 5 (masculino)
$ws.Range("L5").Value = 71.27543548498546
$ws.Range("M5").Value = 40.7941143515294
$ws.Range("N5").Value = 10.77163914429046
$ws.Range("O5").Value = 461.5197709071476
$ws.Range("R5").Value = 214.0671103886876
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = 3474

# Row 6 (outros)
$ws.Range("L6").Value = 45.24320624776205
$ws.Range("M6").Value = 14.93011251059404
$ws.Range("N6").Value = 21.17889830378416
$ws.Range("O6").Value = 63.40051265045815
$ws.Range("R6").Value = 34.72956747704838
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 87
